$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Copy()
$ws.Rows("6:6").Insert()
$ws.Rows("6:6").ClearContents()

$ws.Range("A6:XFD6").Select()
